# simulation_data.xlsx: refresh "Statistics" and "Accidents" sample rows
# with a newer simulation run (2024-08-01 12:37:xx), trim Statistics down
# to 9 data rows, and grow Accidents out to 15 data rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Statistics": rows 2-10 get new values, rows 11-17 are removed.
# ---------------------------------------------------------------------
$stats = $wb.Worksheets.Item("Statistics")

$statsData = @(
    ,("2024-08-01 12:37:40", 94.23897003537142, 8)
    ,("2024-08-01 12:37:42", 94.17275958540921, 16)
    ,("2024-08-01 12:37:44", 75.30477324354547, 21)
    ,("2024-08-01 12:37:46", 64.54746880472618, 29)
    ,("2024-08-01 12:37:48", 51.6774092598811, 35)
    ,("2024-08-01 12:37:50", 43.59545398424419, 35)
    ,("2024-08-01 12:37:52", 36.5299700089092, 36)
    ,("2024-08-01 12:37:54", 28.37919823182034, 38)
    ,("2024-08-01 12:37:56", 22.26924060751913, 37)
)

$r = 2
foreach ($row in $statsData) {
    $stats.Cells.Item($r, 1).Value = $row[0]
    $stats.Cells.Item($r, 2).Value = $row[1]
    $stats.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Drop the now-stale rows 11-17 (A1:C17 -> A1:C10).
$stats.Range("A11:C17").Clear()

# ---------------------------------------------------------------------
# Sheet "Accidents": rows 2-3 get new values, rows 4-16 are added.
# ---------------------------------------------------------------------
$acc = $wb.Worksheets.Item("Accidents")

$accData = @(
    ,("2024-08-01 12:37:44", "Car and Car", "61.98 and 88.71")
    ,("2024-08-01 12:37:44", "Car and Car", "91.53 and 0.00")
    ,("2024-08-01 12:37:46", "Car and Truck", "0.00 and 85.79")
    ,("2024-08-01 12:37:46", "Car and Truck", "78.59 and 76.66")
    ,("2024-08-01 12:37:47", "Car and Car", "88.74 and 0.00")
    ,("2024-08-01 12:37:48", "Car and Truck", "11.65 and 0.00")
    ,("2024-08-01 12:37:48", "Car and Car", "50.96 and 71.65")
    ,("2024-08-01 12:37:48", "Car and Car", "42.46 and 0.00")
    ,("2024-08-01 12:37:49", "Car and Car", "56.47 and 0.00")
    ,("2024-08-01 12:37:50", "Car and Car", "54.49 and 88.75")
    ,("2024-08-01 12:37:52", "Car and Car", "0.00 and 95.16")
    ,("2024-08-01 12:37:54", "Car and Car", "0.00 and 37.35")
    ,("2024-08-01 12:37:54", "Car and Car", "0.00 and 46.56")
    ,("2024-08-01 12:37:55", "Car and Car", "7.27 and 0.00")
    ,("2024-08-01 12:37:57", "Car and Car", "54.27 and 0.00")
)

$r = 2
foreach ($row in $accData) {
    $acc.Cells.Item($r, 1).Value = $row[0]
    $acc.Cells.Item($r, 2).Value = $row[1]
    $acc.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

Write-Host "done"
